$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.539.37"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "3.148.23"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.425"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.18%  "
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D13").Value = "3.689.69"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  +4.98%  "
$ws.Range("D16").Value = "58.571.64"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.142.27"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.516"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.59%  "
$ws.Range("D29").Value = "0.0₃0858"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("E37").Value = "  +7.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("D39").Value = "2.656.89"
$ws.Range("E39").Value = "  +11.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.709"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("E44").Value = "  +4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "3.191.18"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  +12.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.979"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
